# DockPoints.xlsx update
# - New docking path rows (Startpunkt / Zwischenwert / Einfahrt / Endpunkt / 10cm left shift)
# - New shared strings for the new labels
# - Selection moved to F26

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: Startpunkt shifted by the row 14 offset
$ws.Range("B16").Value = "Startpunkt"
$ws.Range("C16").Formula = "=C5+C`$14"
$ws.Range("D16").Formula = "=D5+D`$14"

# Row 17: Zwischenwert shifted by the row 14 offset
$ws.Range("B17").Value = "Zwischenwert"
$ws.Range("C17").Formula = "=C6+C`$14"
$ws.Range("D17").Formula = "=D6+D`$14"

# Row 18: Einfahrt shifted by the row 14 offset
$ws.Range("B18").Value = "Einfahrt"
$ws.Range("C18").Formula = "=C7+C`$14"
$ws.Range("D18").Formula = "=D7+D`$14"

# Row 19: Endpunkt shifted by the row 14 offset
$ws.Range("B19").Value = "Endpunkt"
$ws.Range("C19").Formula = "=C8+C`$14"
$ws.Range("D19").Formula = "=D8+D`$14"

# Row 14: shift vector (10 cm to the left) based on the normalized direction in row 12
$ws.Range("B14").Value = "Verschiebung nach links 10 cm"
$ws.Range("C14").Formula = "=-D12*E14"
$ws.Range("D14").Formula = "=C12*E14"
$ws.Range("E14").Value = 0.1

# Update the saved selection to match the author's final cursor position
$ws.Range("F26").Select()
